$wb = $excel.ActiveWorkbook

# --- "converting" sheet: add the new test case rows ------------------------
$ws = $wb.Worksheets.Item("converting")

# Header / tester-name cell (row 1, merged A1:F1)
$ws.Cells.Item(1, 1).Value = "Tester name:  Manav Dineshbhai Dhameliya"

# Row 3 - Case 1
$ws.Cells.Item(3, 1).Value = "Case 1"
$ws.Cells.Item(3, 2).Value = "123gg"
$ws.Cells.Item(3, 3).Value = 123
$ws.Cells.Item(3, 4).Value = 123
$ws.Cells.Item(3, 5).Value = "Pass"
$ws.Cells.Item(3, 5).HorizontalAlignment = -4152
$ws.Cells.Item(3, 6).Value = "It Should extract the Integer from the input entered"

# Row 4 - Case 2
$ws.Cells.Item(4, 1).Value = "Case 2"
$ws.Cells.Item(4, 2).Value = "fgd"
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = "Pass"
$ws.Cells.Item(4, 5).HorizontalAlignment = -4152
$ws.Cells.Item(4, 6).Value = "it should give 0 a there is no int in the input entered"

# Row 5 - Case 3
$ws.Cells.Item(5, 1).Value = "Case 3"
$ws.Cells.Item(5, 2).Value = "123g32"
$ws.Cells.Item(5, 3).Value = 123
$ws.Cells.Item(5, 4).Value = 12332
$ws.Cells.Item(5, 5).Value = "Fail"
$ws.Cells.Item(5, 5).HorizontalAlignment = -4152
$ws.Cells.Item(5, 6).Value = "This should give all the integer in the input entered"

# Row 6 - Case 4
$ws.Cells.Item(6, 1).Value = "Case 4"
$ws.Cells.Item(6, 2).Value = '" "'
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = "Pass"
$ws.Cells.Item(6, 5).HorizontalAlignment = -4152
$ws.Cells.Item(6, 6).Value = "If no string is entered then it should give 0 as output"

# Row 7 - Case 5
$ws.Cells.Item(7, 1).Value = "Case 5"
$ws.Cells.Item(7, 2).Value = "manav@1246"
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 1246
$ws.Cells.Item(7, 5).Value = "Fail"
$ws.Cells.Item(7, 5).HorizontalAlignment = -4152
$ws.Cells.Item(7, 6).Value = "To see if any special case are entered will it give any output"

# Row 8 - Case 6 (value needs a leading apostrophe so Excel stores it as
# quote-prefixed text, matching the quotePrefix style in the target file)
$ws.Cells.Item(8, 1).Value = "Case 6"
$ws.Cells.Item(8, 2).Value = "'+12@'"
$ws.Cells.Item(8, 3).Value = 12
$ws.Cells.Item(8, 4).Value = 12
$ws.Cells.Item(8, 5).Value = "Pass"
$ws.Cells.Item(8, 5).HorizontalAlignment = -4152
$ws.Cells.Item(8, 6).Value = "This should give all the integer in the input entered"

# Row 9 - Case 7 (same quote-prefix situation)
$ws.Cells.Item(9, 1).Value = "Case 7"
$ws.Cells.Item(9, 2).Value = "'@122eee'"
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 122
$ws.Cells.Item(9, 5).Value = "Fail"
$ws.Cells.Item(9, 5).HorizontalAlignment = -4152
$ws.Cells.Item(9, 6).Value = "It Should extract the Integer from the input entered"

# Row 10 - Case 8
$ws.Cells.Item(10, 1).Value = "Case 8"
$ws.Cells.Item(10, 2).Value = "q"
$ws.Cells.Item(10, 3).Value = "Programs Ends"
$ws.Cells.Item(10, 4).Value = "Program end"
$ws.Cells.Item(10, 5).Value = "Pass"
$ws.Cells.Item(10, 5).HorizontalAlignment = -4152
$ws.Cells.Item(10, 6).Value = "This should give all the integer in the input entered"

# The new data makes "converting" the active/selected sheet & range;
# selecting here also clears tabSelected from whichever sheet had it before
# ("manipulating") and makes "converting" tab 0 the active one.
$ws.Range("A1:F10").Select()
